# Commit: Fri, Jun 05, 2020  3:06:09 PM
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from {DB318A30-BC1D-477C-B780-573630E77BDE} to
#    {BC2F7ABD-8493-49E2-96B4-DEA314109809}.
# 2) The deck's colour theme is switched from the "Integral / Red Violet"
#    palette over to the standard "Office" palette.

$p = $ppt.ActivePresentation

function Set-TableStyle($slideIndex, $styleId) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($styleId)
        }
    }
}

$newTableStyle = "{BC2F7ABD-8493-49E2-96B4-DEA314109809}"
Set-TableStyle 14 $newTableStyle
Set-TableStyle 15 $newTableStyle
Set-TableStyle 16 $newTableStyle

# Re-colour the presentation theme to the standard Office palette.
# Order matches MsoThemeColorSchemeIndex 1-12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeTheme[$i - 1]
}
